$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 38461984
$ws.Range("I28").Value = 38461984
$ws.Range("K28").Value = 38461984
$ws.Range("M28").Value = -38461499

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4500
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1153
$ws.Range("J80").Value = 1124.8572
$ws.Range("L80").Value = 3374.5716
$ws.Range("N80").Value = -5370.571599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1153
$ws.Range("J83").Value = 1124.8572
$ws.Range("L83").Value = 10123.7148
$ws.Range("N83").Value = -20107.7148

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4372.222
$ws.Range("I116").Value = 3220
$ws.Range("J116").Value = 5812.5
$ws.Range("K116").Value = 3220
$ws.Range("L116").Value = 5812.5
$ws.Range("M116").Value = 222
$ws.Range("N116").Value = -12696.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1019.7568
$ws.Range("I135").Value = 826.875
$ws.Range("J135").Value = 2254.2
$ws.Range("K135").Value = 7441.875
$ws.Range("L135").Value = 20287.8
$ws.Range("M135").Value = -4906.875
$ws.Range("N135").Value = -25357.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 12002
$ws.Range("I141").Value = 13002.737
$ws.Range("K141").Value = 39008.211
$ws.Range("M141").Value = -33828.211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5072.618
$ws.Range("I32").Value = 5072.618
$ws.Range("K32").Value = 5072.618
$ws.Range("M32").Value = -4785.618

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 21856.572
$ws.Range("J55").Value = 24999
$ws.Range("L55").Value = 24999
$ws.Range("N55").Value = -25629

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1435.8
$ws.Range("J88").Value = 1456.7858
$ws.Range("L88").Value = 1456.7858
$ws.Range("N88").Value = -2268.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1435.8
$ws.Range("J91").Value = 1456.7858
$ws.Range("L91").Value = 1456.7858
$ws.Range("N91").Value = -4264.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 17701.143
$ws.Range("I122").Value = 1619.9286
$ws.Range("K122").Value = 4859.7858
$ws.Range("M122").Value = -2409.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2177.1667
$ws.Range("I132").Value = 2119.652
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 6358.956
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -3828.956
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5212.091
$ws.Range("I107").Value = 3833.9375
$ws.Range("K107").Value = 3833.9375
$ws.Range("M107").Value = -1913.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5987.077
$ws.Range("I134").Value = 5987.077
$ws.Range("K134").Value = 17961.231
$ws.Range("M134").Value = -15426.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 229.75
$ws.Range("I11").Value = 210
$ws.Range("J11").Value = 249.5
$ws.Range("K11").Value = 210
$ws.Range("L11").Value = 249.5
$ws.Range("M11").Value = -70
$ws.Range("N11").Value = -529.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 68250
$ws.Range("J74").Value = 68250
$ws.Range("L74").Value = 68250
$ws.Range("N74").Value = -69998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 68250
$ws.Range("J77").Value = 68250
$ws.Range("L77").Value = 204750
$ws.Range("N77").Value = -213486

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3247.1765
$ws.Range("I122").Value = 2636.1667
$ws.Range("K122").Value = 7908.500100000001
$ws.Range("M122").Value = -5458.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 41196.168
$ws.Range("I132").Value = 30000
$ws.Range("J132").Value = 43435.4
$ws.Range("K132").Value = 90000
$ws.Range("L132").Value = 130306.2
$ws.Range("M132").Value = -87470
$ws.Range("N132").Value = -135366.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4156.1665
$ws.Range("I134").Value = 3920.8
$ws.Range("J134").Value = 5333
$ws.Range("K134").Value = 11762.4
$ws.Range("L134").Value = 15999
$ws.Range("M134").Value = -9227.400000000001
$ws.Range("N134").Value = -21069

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 37037256
$ws.Range("I18").Value = 250.14285
$ws.Range("K18").Value = 750.4285500000001
$ws.Range("M18").Value = -581.4285500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 919.8570999999999
$ws.Range("J34").Value = 439.08334
$ws.Range("L34").Value = 1317.25002
$ws.Range("N34").Value = -1485.25002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6599.077
$ws.Range("J39").Value = 6599.077
$ws.Range("L39").Value = 19797.231
$ws.Range("N39").Value = -20385.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1442.5625
$ws.Range("J55").Value = 1619.2565
$ws.Range("L55").Value = 4857.7695
$ws.Range("N55").Value = -5211.7695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1298.2222
$ws.Range("I97").Value = 981
$ws.Range("J97").Value = 1932.6666
$ws.Range("K97").Value = 2943
$ws.Range("L97").Value = 5797.9998
$ws.Range("M97").Value = -2447
$ws.Range("N97").Value = -6789.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 14073.444
$ws.Range("I129").Value = 832.625
$ws.Range("J129").Value = 120000
$ws.Range("K129").Value = 2497.875
$ws.Range("L129").Value = 360000
$ws.Range("M129").Value = 2502.125
$ws.Range("N129").Value = -370000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7484.8
$ws.Range("I137").Value = 6856
$ws.Range("K137").Value = 20568
$ws.Range("M137").Value = -15468

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3354000
$ws.Range("I11").Value = 3354000
$ws.Range("K11").Value = 3354000
$ws.Range("M11").Value = -3353861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6142.5713
$ws.Range("I70").Value = 5399.8
$ws.Range("J70").Value = 7999.5
$ws.Range("K70").Value = 5399.8
$ws.Range("L70").Value = 7999.5
$ws.Range("M70").Value = -5129.8
$ws.Range("N70").Value = -8539.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6142.5713
$ws.Range("I73").Value = 5399.8
$ws.Range("J73").Value = 7999.5
$ws.Range("K73").Value = 5399.8
$ws.Range("L73").Value = 7999.5
$ws.Range("M73").Value = -4463.8
$ws.Range("N73").Value = -9871.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2545.7778
$ws.Range("I122").Value = 2126.9614
$ws.Range("K122").Value = 6380.8842
$ws.Range("M122").Value = -3930.8842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4047.9285
$ws.Range("I132").Value = 1882.6666
$ws.Range("K132").Value = 5647.9998
$ws.Range("M132").Value = -3117.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 709
$ws.Range("I22").Value = 699.5
$ws.Range("J22").Value = 728
$ws.Range("K22").Value = 699.5
$ws.Range("L22").Value = 728
$ws.Range("M22").Value = -404.5
$ws.Range("N22").Value = -1318

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 709
$ws.Range("I27").Value = 699.5
$ws.Range("J27").Value = 728
$ws.Range("K27").Value = 699.5
$ws.Range("L27").Value = 728
$ws.Range("M27").Value = -592.5
$ws.Range("N27").Value = -942

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 69000
$ws.Range("J59").Value = 69000
$ws.Range("L59").Value = 69000
$ws.Range("N59").Value = -70308

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5682.6665
$ws.Range("I132").Value = 5819.4
$ws.Range("K132").Value = 17458.2
$ws.Range("M132").Value = -14928.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3815.3667
$ws.Range("I136").Value = 2911.5715
$ws.Range("J136").Value = 4090.4348
$ws.Range("K136").Value = 8734.7145
$ws.Range("L136").Value = 12271.3044
$ws.Range("M136").Value = -6184.7145
$ws.Range("N136").Value = -17371.3044

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2
$ws.Range("N5").Value = -226
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 161709.75
$ws.Range("J16").Value = 175613
$ws.Range("L16").Value = 175613
$ws.Range("N16").Value = -176197

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 798.86664
$ws.Range("I113").Value = 662.2857
$ws.Range("K113").Value = 1986.8571
$ws.Range("M113").Value = 183.1428999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3987.926
$ws.Range("I132").Value = 3979
$ws.Range("J132").Value = 4003.1
$ws.Range("K132").Value = 11937
$ws.Range("L132").Value = 12009.3
$ws.Range("M132").Value = -9407
$ws.Range("N132").Value = -17069.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7218
$ws.Range("I136").Value = 7218
$ws.Range("K136").Value = 21654
$ws.Range("M136").Value = -19104
